$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2; B = "Bitcoin"; C = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"; D = "63.085.80"; E = "  -0.68%  " },
    @{ Row = 3; B = "Ethereum"; C = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"; D = "2.551.42"; E = "  +0.29%  " },
    @{ Row = 4; B = "TetherUSD"; C = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"; D = "'1.00"; E = "  -0.02%  " },
    @{ Row = 5; B = "BNB"; C = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"; D = "'585.65"; E = "  +2.44%  " },
    @{ Row = 6; B = "Solana"; C = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D = "'147.29"; E = "  -2.40%  " },
    @{ Row = 7; B = "USDC"; C = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"; D = "'1.00"; E = "  +0.00%  " },
    @{ Row = 8; B = "XRP"; C = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"; D = "'0.585"; E = "  -0.54%  " },
    @{ Row = 9; B = "Dogecoin"; C = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D = "'0.106"; E = "  -0.49%  " },
    @{ Row = 10; B = "Toncoin"; C = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D = "'5.57"; E = "  -3.11%  " },
    @{ Row = 11; B = "TRON"; C = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D = "'0.152"; E = "  -0.25%  " },
    @{ Row = 12; B = "Cardano"; C = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; D = "'0.355"; E = "  -1.15%  " },
    @{ Row = 13; B = "Avalanche"; C = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D = "'27.51"; E = "  -3.14%  " },
    @{ Row = 14; B = "WrappedliquidstakedEther2.0"; C = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D = "3.006.83"; E = "  +0.21%  " },
    @{ Row = 15; B = "WrappedBTC"; C = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D = "62.997.03"; E = "  -0.67%  " },
    @{ Row = 16; B = "ShibaInu"; C = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D = "'0.0000143"; E = "  -0.57%  " },
    @{ Row = 17; B = "WrappedEther"; C = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D = "2.551.39"; E = "  +0.28%  " },
    @{ Row = 18; B = "Chainlink"; C = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D = "'11.37"; E = "  -2.69%  " },
    @{ Row = 19; B = "BitcoinCash"; C = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D = "'337.73"; E = "  -0.77%  " },
    @{ Row = 20; B = "Polkadot"; C = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D = "'4.34"; E = "  -0.66%  " },
    @{ Row = 21; B = "Uniswap"; C = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D = "'6.77"; E = "  -1.13%  " },
    @{ Row = 22; B = "Dai"; C = "https://coinranking.com/coin/MoTuySvg7+dai-dai"; D = "'1.00"; E = "  -0.01%  " },
    @{ Row = 23; B = "Litecoin"; C = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D = "'65.68"; E = "  -0.82%  " },
    @{ Row = 24; B = "WrappedeETH"; C = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"; D = "2.681.20"; E = "  +0.53%  " },
    @{ Row = 25; B = "Kaspa"; C = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; D = "'0.169"; E = "  -0.50%  " },
    @{ Row = 26; B = "Fetch.AI"; C = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"; D = "'1.61"; E = "  +0.23%  " },
    @{ Row = 27; B = "SuiNetwork"; C = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"; D = "'1.49"; E = "  -0.92%  " },
    @{ Row = 28; B = "Binance-PegBSC-USD"; C = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"; D = "'0.999"; E = "  -0.07%  " },
    @{ Row = 29; B = "InternetComputer(DFINITY)"; C = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D = "'8.37"; E = "  -2.55%  " },
    @{ Row = 30; B = "Aptos"; C = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D = "'7.71"; E = "  +7.63%  " },
    @{ Row = 31; B = "PancakeSwap"; C = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D = "'1.98"; E = "  +5.58%  " },
    @{ Row = 32; B = "PEPE"; C = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"; D = "0.0₃0816"; E = "  -2.06%  " },
    @{ Row = 33; B = "Monero"; C = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D = "'178.17"; E = "  +0.00%  " },
    @{ Row = 34; B = "Bittensor"; C = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"; D = "'419.83"; E = "  -0.27%  " },
    @{ Row = 35; B = "ImmutableX"; C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D = "'1.55"; E = "  -1.26%  " },
    @{ Row = 36; B = "PolygonEcosystemToken"; C = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"; D = "'0.401"; E = "  -1.29%  " },
    @{ Row = 37; B = "EthereumClassic"; C = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D = "'19.15"; E = "  -0.43%  " },
    @{ Row = 38; B = "USDe"; C = "https://coinranking.com/coin/exbfr2U-0+usde-usde"; D = "'0.999"; E = "  +0.03%  " },
    @{ Row = 39; B = "NEARProtocol"; C = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D = "'4.37"; E = "  -2.82%  " },
    @{ Row = 40; B = "Stacks"; C = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"; D = "'1.75"; E = "  -1.11%  " },
    @{ Row = 41; B = "FirstDigitalUSD"; C = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"; D = "'1.00"; E = "  +0.00%  " },
    @{ Row = 42; B = "OKB"; C = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; D = "'39.76"; E = "  -0.17%  " },
    @{ Row = 43; B = "Aave"; C = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"; D = "'150.44"; E = "  -2.41%  " },
    @{ Row = 44; B = "Filecoin"; C = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D = "'3.78"; E = "  -0.68%  " },
    @{ Row = 45; B = "InjectiveProtocol"; C = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; D = "'20.83"; E = "  -1.46%  " },
    @{ Row = 46; B = "Hedera"; C = "https://coinranking.com/coin/jad286TjB+hedera-hbar"; D = "'0.0541"; E = "  +1.71%  " },
    @{ Row = 47; B = "Mantle"; C = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"; D = "'0.603"; E = "  -1.51%  " },
    @{ Row = 48; B = "Stellar"; C = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D = "'0.0972"; E = "  +0.42%  " },
    @{ Row = 49; B = "VeChain"; C = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D = "'0.0237"; E = "  -1.56%  " },
    @{ Row = 50; B = "EnergySwap"; C = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D = "'18.32"; E = "  -1.80%  " },
    @{ Row = 51; B = "dogwifhat"; C = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"; D = "'1.72"; E = "  -5.65%  " }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}
